# TST: Test for float --> int in the middle of object dtype
#
# pandas/io/tests/data/test_types.xlsx — Str2Col originally held the
# strings a/b/c/d/e in E2:E6. Replace the "b" in E3 with the numeric
# value 3, so the column mixes floats/ints with strings in the middle
# (object dtype edge case). Excel drops the now-unreferenced "b" shared
# string automatically and reindexes the shared-string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = 3

# Matches the selection/window state recorded in the saved workbook.
$ws.Range("E4").Select()

$win = $wb.Windows.Item(1)
$win.Left = 6200
$win.Top = 2220
